$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @(
  @("ECs", "Il17c", "Il17re", "ECs", 3, 1, 2.742348, 8.227043999999999, 0.8112951562181048, 0.8112951562181048, 2, 0.6666666666666666, 0.001470666666666667, 0.004412, 0.012032749150453, 0.012032749150453, 0.004033079792, 0.036297718128, 0.009762111101750036, 0.009762111101750034),
  @("ECs", "Il17c", "Il17re", "FAPs", 3, 1, 2.742348, 8.227043999999999, 0.8112951562181048, 0.8112951562181048, 1, 0.3333333333333333, 0.05659733333333333, 0.169792, 0.4630699328544234, 0.4630699328544233, 0.155209583872, 1.396886254848, 0.3756863935150367, 0.3756863935150367),
  @("ECs", "Il17c", "Il17re", "sCs", 3, 1, 2.742348, 8.227043999999999, 0.8112951562181048, 0.8112951562181048, 3, 1, 0.064154, 0.192462, 0.5248973179951236, 0.5248973179951236, 0.175932593592, 1.583393342328, 0.425846651601318, 0.425846651601318),
  @("FAPs", "Il17c", "Il17re", "ECs", 3, 1, 0.3936376666666666, 1.180913, 0.1164536128425946, 0.1164536128425946, 2, 0.6666666666666666, 0.001470666666666667, 0.004412, 0.012032749150453, 0.012032749150453, 0.000578909795111111, 0.005210188156, 0.001401257110998913, 0.001401257110998912),
  @("FAPs", "Il17c", "Il17re", "FAPs", 3, 1, 0.3936376666666666, 1.180913, 0.1164536128425946, 0.1164536128425946, 1, 0.3333333333333333, 0.05659733333333333, 0.169792, 0.4630699328544234, 0.4630699328544233, 0.02227884223288889, 0.200509580096, 0.05392616667967529, 0.05392616667967529),
  @("FAPs", "Il17c", "Il17re", "sCs", 3, 1, 0.3936376666666666, 1.180913, 0.1164536128425946, 0.1164536128425946, 3, 1, 0.064154, 0.192462, 0.5248973179951236, 0.5248973179951236, 0.02525343086733333, 0.227280877806, 0.06112618905192038, 0.06112618905192038),
  @("sCs", "Il17c", "Il17re", "ECs", 3, 1, 0.2442243333333333, 0.732673, 0.07225123093930062, 0.07225123093930062, 2, 0.6666666666666666, 0.001470666666666667, 0.004412, 0.012032749150453, 0.012032749150453, 0.0003591725862222223, 0.003232553276, 0.0008693809377040532, 0.0008693809377040531),
  @("sCs", "Il17c", "Il17re", "FAPs", 3, 1, 0.2442243333333333, 0.732673, 0.07225123093930062, 0.07225123093930062, 1, 0.3333333333333333, 0.05659733333333333, 0.169792, 0.4630699328544234, 0.4630699328544233, 0.01382244600177778, 0.124402014016, 0.03345737265971137, 0.03345737265971137),
  @("sCs", "Il17c", "Il17re", "sCs", 3, 1, 0.2442243333333333, 0.732673, 0.07225123093930062, 0.07225123093930062, 3, 1, 0.064154, 0.192462, 0.5248973179951236, 0.5248973179951236, 0.01566796788066667, 0.141011710926, 0.03792447734188519, 0.03792447734188519)
)

$startRow = 2
for ($i = 0; $i -lt $rowsData.Length; $i++) {
  $row = $rowsData[$i]
  $r = $startRow + $i
  for ($c = 0; $c -lt $row.Length; $c++) {
    $ws.Cells.Item($r, $c + 1).Value = $row[$c]
  }
}

